$d = $word.ActiveDocument

# --- 1. Collapse the three CORE COMPETENCIES detail paragraphs into one summary line ---
# Paragraph 6 holds "Product Management & Strategy: ..." (long form) -- keep this
# paragraph but replace its text, then delete the following two long paragraphs
# (Technical Product Development / Platform & Infrastructure details).
$core1 = $d.Paragraphs(6)
$core1.Range.Text = "Product Management & Strategy " + [char]0x2022 + " Technical Product Development " + [char]0x2022 + " Platform & Infrastructure"

# After editing paragraph 6's text in place, paragraphs 7 and 8 are still the old
# long-form "Technical Product Development" and "Platform & Infrastructure" paragraphs.
$d.Paragraphs(7).Range.Delete()
$d.Paragraphs(7).Range.Delete()

# --- 2. Append a new "TECHNICAL SKILLS" section at the end of the document ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs($d.Paragraphs.Count)
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading2"

$headingPara.Range.InsertParagraphAfter()
$skill1 = $d.Paragraphs($d.Paragraphs.Count)
$skill1.Style = "Normal"
$skill1.Range.Text = "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics"

$skill1.Range.InsertParagraphAfter()
$skill2 = $d.Paragraphs($d.Paragraphs.Count)
$skill2.Style = "Normal"
$skill2.Range.Text = "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration"

$skill2.Range.InsertParagraphAfter()
$skill3 = $d.Paragraphs($d.Paragraphs.Count)
$skill3.Style = "Normal"
$skill3.Range.Text = "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training"
